$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 28894.18
$ws.Range("J87").Value = 28894.18
$ws.Range("L87").Value = 28894.18
$ws.Range("N87").Value = -31390.18
$ws.Range("H90").Value = 28894.18
$ws.Range("J90").Value = 28894.18
$ws.Range("L90").Value = 86682.54000000001
$ws.Range("N90").Value = -99162.54000000001
$ws.Range("H125").Value = 1751.6666
$ws.Range("I125").Value = 1416
$ws.Range("J125").Value = 1847.5714
$ws.Range("K125").Value = 12744
$ws.Range("L125").Value = 16628.1426
$ws.Range("M125").Value = -10284
$ws.Range("N125").Value = -21548.1426
$ws.Range("H135").Value = 26162.25
$ws.Range("I135").Value = 31266.908
$ws.Range("J135").Value = 2097.4285
$ws.Range("K135").Value = 281402.172
$ws.Range("L135").Value = 18876.8565
$ws.Range("M135").Value = -278867.172
$ws.Range("N135").Value = -23946.8565
$ws.Range("H138").Value = 2262.0488
$ws.Range("I138").Value = 1426.0834
$ws.Range("J138").Value = 3442.2354
$ws.Range("K138").Value = 4278.2502
$ws.Range("L138").Value = 10326.7062
$ws.Range("M138").Value = 861.7497999999996
$ws.Range("N138").Value = -20606.7062

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 30355
$ws.Range("J24").Value = 30355
$ws.Range("L24").Value = 30355
$ws.Range("N24").Value = -31103
$ws.Range("H32").Value = 18158.285
$ws.Range("I32").Value = 3280.5916
$ws.Range("J32").Value = 99413.38
$ws.Range("K32").Value = 3280.5916
$ws.Range("L32").Value = 99413.38
$ws.Range("M32").Value = -2993.5916
$ws.Range("N32").Value = -99987.38
$ws.Range("H100").Value = 30355
$ws.Range("J100").Value = 30355
$ws.Range("L100").Value = 30355
$ws.Range("N100").Value = -32519
$ws.Range("H132").Value = 296964.47
$ws.Range("I132").Value = 456942.38
$ws.Range("J132").Value = 3671.6667
$ws.Range("K132").Value = 1370827.14
$ws.Range("L132").Value = 11015.0001
$ws.Range("M132").Value = -1368297.14
$ws.Range("N132").Value = -16075.0001

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 92618.27
$ws.Range("I134").Value = 108728.68
$ws.Range("J134").Value = 2400
$ws.Range("K134").Value = 326186.04
$ws.Range("L134").Value = 7200
$ws.Range("M134").Value = -323651.04
$ws.Range("N134").Value = -12270

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 52000
$ws.Range("I23").Value = 52000
$ws.Range("K23").Value = 52000
$ws.Range("M23").Value = -51760
$ws.Range("H27").Value = 52000
$ws.Range("I27").Value = 52000
$ws.Range("K27").Value = 52000
$ws.Range("M27").Value = -51808
$ws.Range("H31").Value = 1610.3846
$ws.Range("I31").Value = 1311.56
$ws.Range("J31").Value = 2144
$ws.Range("K31").Value = 1311.56
$ws.Range("L31").Value = 2144
$ws.Range("M31").Value = -1016.56
$ws.Range("N31").Value = -2734
$ws.Range("H34").Value = 1610.3846
$ws.Range("I34").Value = 1311.56
$ws.Range("J34").Value = 2144
$ws.Range("K34").Value = 1311.56
$ws.Range("L34").Value = 2144
$ws.Range("M34").Value = -1109.56
$ws.Range("N34").Value = -2548
$ws.Range("H96").Value = 21500
$ws.Range("J96").Value = 21500
$ws.Range("L96").Value = 21500
$ws.Range("N96").Value = -26992
$ws.Range("H132").Value = 3548.3333
$ws.Range("I132").Value = 3058.2
$ws.Range("K132").Value = 9174.599999999999
$ws.Range("M132").Value = -6644.599999999999
$ws.Range("H134").Value = 5549
$ws.Range("I134").Value = 6046.1665
$ws.Range("J134").Value = 3162.6
$ws.Range("K134").Value = 18138.4995
$ws.Range("L134").Value = 9487.799999999999
$ws.Range("M134").Value = -15603.4995
$ws.Range("N134").Value = -14557.8

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 3930
$ws.Range("I134").Value = 300
$ws.Range("J134").Value = 4333.3335
$ws.Range("K134").Value = 900
$ws.Range("L134").Value = 13000.0005
$ws.Range("M134").Value = 4170
$ws.Range("N134").Value = -23140.0005

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H25").Value = 10000
$ws.Range("J25").Value = 10000
$ws.Range("L25").Value = 10000
$ws.Range("N25").Value = -11058
$ws.Range("H122").Value = 2478.8948
$ws.Range("I122").Value = 2660.7273
$ws.Range("J122").Value = 2228.875
$ws.Range("K122").Value = 7982.1819
$ws.Range("L122").Value = 6686.625
$ws.Range("M122").Value = -5532.1819
$ws.Range("N122").Value = -11586.625
$ws.Range("H126").Value = 3733.487
$ws.Range("I126").Value = 2053.3684
$ws.Range("J126").Value = 5329.6
$ws.Range("K126").Value = 6160.1052
$ws.Range("L126").Value = 15988.8
$ws.Range("M126").Value = -3690.1052
$ws.Range("N126").Value = -20928.8
$ws.Range("H132").Value = 3744.7778
$ws.Range("I132").Value = 3540
$ws.Range("J132").Value = 4000.75
$ws.Range("K132").Value = 10620
$ws.Range("L132").Value = 12002.25
$ws.Range("M132").Value = -8090
$ws.Range("N132").Value = -17062.25

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1283.8
$ws.Range("I22").Value = 1022.2222
$ws.Range("J22").Value = 1430.9375
$ws.Range("K22").Value = 1022.2222
$ws.Range("L22").Value = 1430.9375
$ws.Range("M22").Value = -727.2222
$ws.Range("N22").Value = -2020.9375
$ws.Range("H27").Value = 1283.8
$ws.Range("I27").Value = 1022.2222
$ws.Range("J27").Value = 1430.9375
$ws.Range("K27").Value = 1022.2222
$ws.Range("L27").Value = 1430.9375
$ws.Range("M27").Value = -915.2222
$ws.Range("N27").Value = -1644.9375
$ws.Range("H136").Value = 1544.7241
$ws.Range("I136").Value = 1605
$ws.Range("K136").Value = 4815
$ws.Range("M136").Value = -2265

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H30").Value = 16540
$ws.Range("I30").Value = 8000
$ws.Range("J30").Value = 18675
$ws.Range("K30").Value = 8000
$ws.Range("L30").Value = 18675
$ws.Range("M30").Value = -7893
$ws.Range("N30").Value = -18889
$ws.Range("H101").Value = 9232
$ws.Range("J101").Value = 9232
$ws.Range("L101").Value = 9232
$ws.Range("N101").Value = -15722
$ws.Range("H132").Value = 1873.6279
$ws.Range("I132").Value = 1701.8379
$ws.Range("J132").Value = 2933
$ws.Range("K132").Value = 5105.5137
$ws.Range("L132").Value = 8799
$ws.Range("M132").Value = -2575.5137
$ws.Range("N132").Value = -13859
